$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.4480592728037891
$ws.Range("C2").Value = 0.05689896751667334
$ws.Range("D2").Value = 0.07858837496944204
$ws.Range("E2").Value = 0.4132022044686465
$ws.Range("G2").Value = 0.002439320690758406
$ws.Range("K2").Value = 0.428420881401621
$ws.Range("N2").Value = 1.483451473734611
$ws.Range("O2").Value = 3.098667715785666

$ws.Range("B3").Value = 0.4037757052505526
$ws.Range("C3").Value = 0.05005304256543752
$ws.Range("D3").Value = 0.07127684251783251
$ws.Range("E3").Value = 0.3605282455053498
$ws.Range("G3").Value = 0.00244236072053212
$ws.Range("K3").Value = 0.3799761155893862
$ws.Range("N3").Value = 1.500844144952591
$ws.Range("O3").Value = 3.089321637985989

$ws.Range("B4").Value = 0.3767124667835162
$ws.Range("C4").Value = 0.04583028563760649
$ws.Range("D4").Value = 0.0668227801316732
$ws.Range("E4").Value = 0.3282780526915587
$ws.Range("G4").Value = 0.002444325254653591
$ws.Range("K4").Value = 0.3502968897492735
$ws.Range("N4").Value = 1.512073873888571
$ws.Range("O4").Value = 3.085493650046459

$ws.Range("B5").Value = 0.3657160698450639
$ws.Range("C5").Value = 0.04410460447228104
$ws.Range("D5").Value = 0.06501655246559324
$ws.Range("E5").Value = 0.3151571629653631
$ws.Range("G5").Value = 0.002445150526642396
$ws.Range("K5").Value = 0.3382190719133575
$ws.Range("N5").Value = 1.516788368825747
$ws.Range("O5").Value = 3.084413133244681

$ws.Range("B6").Value = 0.3638920691935539
$ws.Range("C6").Value = 0.04381776174791696
$ws.Range("D6").Value = 0.06471716288500318
$ws.Range("E6").Value = 0.3129796858374618
$ws.Range("G6").Value = 0.002445289057009643
$ws.Range("K6").Value = 0.3362145715712757
$ws.Range("N6").Value = 1.517579555016024
$ws.Range("O6").Value = 3.084262640892462

$ws.Range("B7").Value = 0.3765640353935282
$ws.Range("C7").Value = 0.04580703224610261
$ws.Range("D7").Value = 0.06679838495688273
$ws.Range("E7").Value = 0.3281010157518836
$ws.Range("G7").Value = 0.002444336284352699
$ws.Range("K7").Value = 0.3501339362926785
$ws.Range("N7").Value = 1.512136895497928
$ws.Range("O7").Value = 3.085477138011157

$ws.Range("B8").Value = 0.4327640008776825
$ws.Range("C8").Value = 0.05454249519971199
$ws.Range("D8").Value = 0.07606002570723547
$ws.Range("E8").Value = 0.3950200837645355
$ws.Range("G8").Value = 0.002440348612912386
$ws.Range("K8").Value = 0.4117034584250803
$ws.Range("N8").Value = 1.489334107619165
$ws.Range("O8").Value = 3.095048038384107

$ws.Range("B9").Value = 0.5439796911210806
$ws.Range("C9").Value = 0.07152072046636704
$ws.Range("D9").Value = 0.09450384746519092
$ws.Range("E9").Value = 0.5270626341443148
$ws.Range("G9").Value = 0.002433302291184936
$ws.Range("K9").Value = 0.5329673495053555
$ws.Range("N9").Value = 1.448992527284084
$ws.Range("O9").Value = 3.12902664994391

$ws.Range("B10").Value = 0.626312568815905
$ws.Range("C10").Value = 0.08390509927934886
$ws.Range("D10").Value = 0.1082306247403295
$ws.Range("E10").Value = 0.6246985393689783
$ws.Range("G10").Value = 0.002428591783212716
$ws.Range("K10").Value = 0.6223930653313516
$ws.Range("N10").Value = 1.422027968160325
$ws.Range("O10").Value = 3.163341697906418

$ws.Range("B11").Value = 0.6639055258974906
$ws.Range("C11").Value = 0.08952039028959291
$ws.Range("D11").Value = 0.1145145458197874
$ws.Range("E11").Value = 0.6692794292695652
$ws.Range("G11").Value = 0.002426549039695274
$ws.Range("K11").Value = 0.6631507444783722
$ws.Range("N11").Value = 1.41034345571968
$ws.Range("O11").Value = 3.181000559544287

$ws.Range("B12").Value = 0.6781610370887279
$ws.Range("C12").Value = 0.09164414064673565
$ws.Range("D12").Value = 0.1168998393693528
$ws.Range("E12").Value = 0.6861870194360193
$ws.Range("G12").Value = 0.002425789815751684
$ws.Range("K12").Value = 0.6785958202061977
$ws.Range("N12").Value = 1.406002676137067
$ws.Range("O12").Value = 3.187983419499744

$ws.Range("B13").Value = 0.675089976282095
$ws.Range("C13").Value = 0.09118687017169691
$ws.Range("D13").Value = 0.1163858694477256
$ws.Range("E13").Value = 0.6825444901209181
$ws.Range("G13").Value = 0.002425952692556228
$ws.Range("K13").Value = 0.6752689581234108
$ws.Range("N13").Value = 1.406933805879593
$ws.Range("O13").Value = 3.186466358873474

$ws.Range("B14").Value = 0.6650779379752407
$ws.Range("C14").Value = 0.08969516532290811
$ws.Range("D14").Value = 0.1147106708065451
$ws.Range("E14").Value = 0.6706698992048246
$ws.Range("G14").Value = 0.002426486291168384
$ws.Range("K14").Value = 0.6644211997873981
$ws.Range("N14").Value = 1.409984654527314
$ws.Range("O14").Value = 3.181569107681071

$ws.Range("B15").Value = 0.6589478589540647
$ws.Range("C15").Value = 0.0887811097635165
$ws.Range("D15").Value = 0.1136853074358015
$ws.Range("E15").Value = 0.6633997836101173
$ws.Range("G15").Value = 0.002426814998438354
$ws.Range("K15").Value = 0.6577780673564462
$ws.Range("N15").Value = 1.411864319471892
$ws.Range("O15").Value = 3.178607962454578

$ws.Range("B16").Value = 0.6238585703800652
$ws.Range("C16").Value = 0.08353775697099763
$ws.Range("D16").Value = 0.1078207533576858
$ws.Range("E16").Value = 0.6217885746849277
$ws.Range("G16").Value = 0.002428727289110012
$ws.Range("K16").Value = 0.6197310043177708
$ws.Range("N16").Value = 1.422803305362187
$ws.Range("O16").Value = 3.162228980456632

$ws.Range("B17").Value = 0.6023679964891642
$ws.Range("C17").Value = 0.08031641836734593
$ws.Range("D17").Value = 0.1042331832064036
$ws.Range("E17").Value = 0.5963052844230674
$ws.Range("G17").Value = 0.002429926000061284
$ws.Range("K17").Value = 0.5964101242134632
$ws.Range("N17").Value = 1.429663193938524
$ws.Range("O17").Value = 3.152706610919921

$ws.Range("B18").Value = 0.5900203093962091
$ws.Range("C18").Value = 0.07846185465092503
$ws.Range("D18").Value = 0.102173428174865
$ws.Range("E18").Value = 0.5816634618532817
$ws.Range("G18").Value = 0.002430624892210231
$ws.Range("K18").Value = 0.5830038693969186
$ws.Range("N18").Value = 1.4336635682318
$ws.Range("O18").Value = 3.147422362354519

$ws.Range("B19").Value = 0.5858418534452028
$ws.Range("C19").Value = 0.07783363293845014
$ws.Range("D19").Value = 0.1014766687985258
$ws.Range("E19").Value = 0.5767085945937396
$ws.Range("G19").Value = 0.002430863146245922
$ws.Range("K19").Value = 0.5784659989246563
$ws.Range("N19").Value = 1.43502742231189
$ws.Range("O19").Value = 3.145666277434827

$ws.Range("B20").Value = 0.6046543477949058
$ws.Range("C20").Value = 0.0806595148461895
$ws.Range("D20").Value = 0.104614700957967
$ws.Range("E20").Value = 0.5990164069219617
$ws.Range("G20").Value = 0.00242979742024371
$ws.Range("K20").Value = 0.5988919171689986
$ws.Range("N20").Value = 1.428927279236945
$ws.Range("O20").Value = 3.153700323970554

$ws.Range("B21").Value = 0.6680181781960641
$ws.Range("C21").Value = 0.09013338671294946
$ws.Range("D21").Value = 0.1152025621294968
$ws.Range("E21").Value = 0.6741570394139984
$ws.Range("G21").Value = 0.002426329172390309
$ws.Range("K21").Value = 0.6676071509722306
$ws.Range("N21").Value = 1.409086267679022
$ws.Range("O21").Value = 3.182999510472996

$ws.Range("B22").Value = 0.7095458748505052
$ws.Range("C22").Value = 0.0963097362187284
$ws.Range("D22").Value = 0.1221556287339212
$ws.Range("E22").Value = 0.7234169304066143
$ws.Range("G22").Value = 0.002424145898990014
$ws.Range("K22").Value = 0.712580753773409
$ws.Range("N22").Value = 1.396608114160479
$ws.Range("O22").Value = 3.203873134691946

$ws.Range("B23").Value = 0.6873712204096307
$ws.Range("C23").Value = 0.09301470507438125
$ws.Range("D23").Value = 0.1184415928086935
$ws.Range("E23").Value = 0.6971115273134529
$ws.Range("G23").Value = 0.002425303543885106
$ws.Range("K23").Value = 0.6885716553197483
$ws.Range("N23").Value = 1.403223110442344
$ws.Range("O23").Value = 3.192574252555232

$ws.Range("B24").Value = 0.6036206649034455
$ws.Range("C24").Value = 0.08050440894604094
$ws.Range("D24").Value = 0.1044422081050413
$ws.Range("E24").Value = 0.5977906810258276
$ws.Range("G24").Value = 0.002429855520827943
$ws.Range("K24").Value = 0.5977698947378371
$ws.Range("N24").Value = 1.429259810177591
$ws.Range("O24").Value = 3.153250473619835

$ws.Range("B25").Value = 0.5137841628844058
$ws.Range("C25").Value = 0.06694361101070001
$ws.Range("D25").Value = 0.08948369256359001
$ws.Range("E25").Value = 0.4912412928645296
$ws.Range("G25").Value = 0.002435126231089175
$ws.Range("K25").Value = 0.5001046564064779
$ws.Range("N25").Value = 1.459436900154866
$ws.Range("O25").Value = 3.118198119162344
